# Inventory.xlsx edit: add a "Code" column (E) with part numbers, fix a
# couple of mislabeled items, correct a price, and switch the Total
# formula from additive to multiplicative (Qty * Unit Price).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column E header + per-row part codes -----------------------
# (Writing these in row order first so the shared-string table grows in
#  the same sequence as the source edit.)
$ws.Range("E1").Value  = "Code"
$ws.Range("E2").Value  = "6280K311"
$ws.Range("E3").Value  = "6280K122"
$ws.Range("E4").Value  = "6280K127"
$ws.Range("E5").Value  = "2780T37"
$ws.Range("E6").Value  = "2780T38"
$ws.Range("E7").Value  = "60635K2"
$ws.Range("E8").Value  = "7265K2"
$ws.Range("E9").Value  = "2781T41"
$ws.Range("E10").Value = "5947K62"
$ws.Range("E11").Value = "1886K15"
$ws.Range("E12").Value = "1497K31"

# --- Row 12: relabel item + corrected unit price ----------------------
$ws.Range("A12").Value = "1/2"" Keyed Rotary Shaft"
$ws.Range("D12").Value = 29.37

# --- Remaining new Code values -----------------------------------------
$ws.Range("E13").Value = "**Ranges from Different Prices for different materials "
$ws.Range("E14").Value = "1655T43"
$ws.Range("E15").Value = "1655T918"
$ws.Range("E16").Value = "88685K984"
$ws.Range("E17").Value = "1388K201"

# --- Column E width, matching the other data columns' look -------------
$ws.Columns.Item(5).ColumnWidth = 16.25

# --- Total formula: multiply quantity by unit price per row ------------
$ws.Range("D18").Formula = "=(C2*D2)+(C3*D3)+(C4*D4)+(C5*D5)+(C6*D6)+(C7*D7)+(C8*D8)+(C9*D9)+(C10*D10)+(C11*D11)+(C12*D12)+(C13*D13)+(C14*D14)+(C15*D15)+(C16*D16)+(C17*D17)"
$ws.Range("D18").ClearFormats()

# --- Selection cosmetics (matches where the author last clicked) -------
[void]$ws.Range("E6").Select()
